# planeacion general del proyecto
# Fix the "Requisitos" breakdown row in the Gantt-style planning sheet:
#  - correct the typo "Runcionales" -> "Funcionales"
#  - make the "Requisitos" stage label bold
#  - mark day 3 ("C9") was already marked; also mark day 3 for "Casos de uso y
#    diagramas" (D9) and day 4 for "modelo de objetos" (E10)
#  - remove the stray mark that had been left in "Programación del sistema" (D14)
#  - leave the cursor on A7, matching where the author left it

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo in the "Requisitos" sub-item label
$ws.Range("A6").Value = "Funcionales y no funcionales"

# Make the "Requisitos" stage label bold
$ws.Range("A5").Font.Bold = $true

# Mark D9 ("Casos de uso y diagramas" / Día 3), copying the existing
# marker-cell formatting (Webdings "a" glyph) from a neighboring marked cell
$ws.Range("C9").Copy() | Out-Null
$ws.Range("D9").PasteSpecial(-4122) | Out-Null
$ws.Range("D9").Value = "a"

# Mark E10 ("modelo de objetos" / Día 4) the same way
$ws.Range("D10").Copy() | Out-Null
$ws.Range("E10").PasteSpecial(-4122) | Out-Null
$ws.Range("E10").Value = "a"

# Clear the stray mark in D14 ("Programación del sistema" / Día 3)
$ws.Range("D14").ClearContents()

# Leave selection where the author left it when saving
$ws.Range("A7").Select() | Out-Null

$excel.CutCopyMode = $false
